$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 3500
$ws.Range('F4').Value = 371
$ws.Range('F5').Value = 8207
$ws.Range('F7').Value = 83
$ws.Range('F8').Value = 2158
$ws.Range('F9').Value = 9
$ws.Range('F12').Value = 1185
$ws.Range('F14').Value = 39
$ws.Range('F16').Value = 572
$ws.Range('F17').Value = 72
$ws.Range('F18').Value = 72
$ws.Range('F19').Value = 435
$ws.Range('F21').Value = 7106
$ws.Range('F23').Value = 55004
$ws.Range('F24').Value = 55004
$ws.Range('F25').Value = 4367
$ws.Range('F26').Value = 1030
$ws.Range('F27').Value = 849
$ws.Range('F31').Value = 9
$ws.Range('F33').Value = 2895
$ws.Range('F35').Value = 26
$ws.Range('F36').Value = 19
$ws.Range('F37').Value = 857
$ws.Range('F38').Value = 1170
$ws.Range('F39').Value = 854
$ws.Range('F40').Value = 147
$ws.Range('F43').Value = 698
$ws.Range('F46').Value = 9
$ws.Range('F47').Value = 145
$ws.Range('F49').Value = 36
$ws.Range('F50').Value = 2466
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F2').Value = 55
$ws.Range('F12').Value = 105
$ws.Range('F16').Value = 7447
$ws.Range('G16').Value = 680
$ws.Range('F17').Value = 103
$ws.Range('F27').Value = 5
$ws.Range('F31').Value = 77
$ws.Range('F32').Value = 20
$ws.Range('F35').Value = 40
$ws.Range('F37').Value = 3
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 2278
$ws.Range('F5').Value = 1537
$ws.Range('F7').Value = 654
$ws.Range('F9').Value = 9323
$ws.Range('G9').Value = '暂时售罄'
$ws.Range('F10').Value = 1641
$ws.Range('F15').Value = 151
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 371
$ws.Range('F4').Value = 8207
$ws.Range('F5').Value = 654
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = '2024-08-17'
$ws.Range('C6').Value = '上海·大悦城·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅'
$ws.Range('D6').Value = '西藏北路166静安大悦城北座6楼611号 次元波板糖'
$ws.Range('E6').Value = '2024.08.17 00:00-10.27 23:59'
$ws.Range('F6').Value = 9323
$ws.Range('H6').Value = 'https://show.bilibili.com/platform/detail.html?id=90438'
$ws.Range('I6').Value = '//i2.hdslb.com/bfs/openplatform/202408/qUE9n4UR1723020534077.png'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = '2024-09-06'
$ws.Range('C7').Value = '上海·「HUNTER×HUNTER × animate cafe」'
$ws.Range('D7').Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws.Range('E7').Value = '2024.09.06 00:00-10.08 23:59'
$ws.Range('F7').Value = 1641
$ws.Range('G7').Value = 30
$ws.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=91069'
$ws.Range('I7').Value = '//i2.hdslb.com/bfs/openplatform/202408/4GkLI2cn1724227065219.jpeg'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = '2024-09-09'
$ws.Range('C8').Value = '上海·日漫咖啡体验'
$ws.Range('E8').Value = '2024.09.09 10:00-12.31 22:00'
$ws.Range('F8').Value = 77
$ws.Range('G8').Value = 60
$ws.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=91993'
$ws.Range('I8').Value = '//i2.hdslb.com/bfs/openplatform/202409/IV5rInWT1725347808557.jpeg'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = '2024-09-10'
$ws.Range('C9').Value = '上海·迷你四驱车赛场'
$ws.Range('D9').Value = '虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶'
$ws.Range('E9').Value = '2024.09.10 10:00-12.31 22:00'
$ws.Range('F9').Value = 4
$ws.Range('G9').Value = 48
$ws.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=92042'
$ws.Range('I9').Value = '//i2.hdslb.com/bfs/openplatform/202409/LzFT5TMO1725348229429.png'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = '2024-09-14'
$ws.Range('C10').Value = '上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）'
$ws.Range('D10').Value = '龙台路10号2F 上海国际传媒港艺术中心'
$ws.Range('E10').Value = '2024.09.14 10:00-10.31 20:00'
$ws.Range('F10').Value = 83
$ws.Range('G10').Value = 88
$ws.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=91856'
$ws.Range('I10').Value = '//i0.hdslb.com/bfs/openplatform/202409/wSR0yFfg1725432304586.jpeg'
$ws.Range('F12').Value = 9
$ws.Range('F16').Value = 39
$ws.Range('F17').Value = 72
$ws.Range('F18').Value = 72
$ws.Range('F19').Value = 435
$ws.Range('F20').Value = 55006
$ws.Range('F22').Value = 1030
$ws.Range('F23').Value = 849
$ws.Range('F27').Value = 105
$ws.Range('F30').Value = 26
$ws.Range('F31').Value = 19
$ws.Range('F32').Value = 1170
$ws.Range('F33').Value = 103
$ws.Range('F35').Value = 147
$ws.Range('F37').Value = 698
$ws.Range('F41').Value = 5
$ws.Range('F44').Value = 145
$ws.Range('F46').Value = 36
$ws.Range('F47').Value = 40
$ws.Range('F49').Value = 2466